$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set a cell value while guaranteeing it stays text even when it looks numeric
# (e.g. "1.000", "6.400") - matches how the crypto price strings are stored in the sheet.
function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        $cell.NumberFormat = '@'
        $cell.Value = $value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $value
    }
}

# Updated Coin / Link / Price / Volume(1h) cells (ranking refresh)
Set-TextValue 'D2' '22.475.61'
$ws.Range('E2').Value = '  +0.39%  '
Set-TextValue 'D3' '1.575.38'
$ws.Range('E3').Value = '  +0.95%  '
Set-TextValue 'D6' '288.14'
$ws.Range('E6').Value = '  +0.78%  '
Set-TextValue 'D7' '0.3701'
$ws.Range('E7').Value = '  +1.78%  '
Set-TextValue 'D8' '47.81'
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('E10').Value = '  +2.44%  '
Set-TextValue 'D11' '0.07573'
$ws.Range('E11').Value = '  +2.60%  '
Set-TextValue 'D12' '1.000'
$ws.Range('E12').Value = '  -0.09%  '
Set-TextValue 'D13' '20.86'
$ws.Range('E13').Value = '  +0.61%  '
Set-TextValue 'D14' '5.961'
$ws.Range('E14').Value = '  +0.87%  '
Set-TextValue 'D15' '6.952'
$ws.Range('E15').Value = '  +1.31%  '
Set-TextValue 'D16' '1.567.34'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('E17').Value = '  +2.24%  '
Set-TextValue 'D18' '88.32'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D20' '6.400'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D21' '0.9997'
$ws.Range('E21').Value = '  -0.13%  '
Set-TextValue 'D22' '16.56'
$ws.Range('E22').Value = '  +3.57%  '
$ws.Range('E23').Value = '  +1.28%  '
Set-TextValue 'D24' '22.471.75'
$ws.Range('E24').Value = '  +0.44%  '
Set-TextValue 'D25' '2.387'
$ws.Range('E25').Value = '  +0.09%  '
Set-TextValue 'D26' '2.639'
$ws.Range('E26').Value = '  +3.71%  '
Set-TextValue 'D27' '151.28'
$ws.Range('E27').Value = '  +1.45%  '
Set-TextValue 'D28' '19.69'
$ws.Range('E28').Value = '  +1.72%  '
Set-TextValue 'D29' '4.995'
$ws.Range('E29').Value = '  -0.34%  '
Set-TextValue 'D30' '125.54'
$ws.Range('E30').Value = '  +2.23%  '
Set-TextValue 'D31' '1.746.46'
$ws.Range('E31').Value = '  +0.65%  '
Set-TextValue 'D32' '1.094'
$ws.Range('E32').Value = '  +4.12%  '
Set-TextValue 'D33' '6.120'
$ws.Range('E33').Value = '  +0.51%  '
Set-TextValue 'D34' '1.987'
$ws.Range('E34').Value = '  -0.02%  '
Set-TextValue 'D35' '9.890'
$ws.Range('E35').Value = '  +3.53%  '
Set-TextValue 'D36' '0.08360'
$ws.Range('E36').Value = '  +1.65%  '
Set-TextValue 'D37' '0.02465'
$ws.Range('E37').Value = '  +4.17%  '
Set-TextValue 'D38' '0.2242'
$ws.Range('E38').Value = '  +1.65%  '
Set-TextValue 'D39' '0.06394'
$ws.Range('E39').Value = '  +0.94%  '
$ws.Range('E40').Value = '  +0.41%  '
Set-TextValue 'D41' '5.370'
$ws.Range('E41').Value = '  +1.36%  '
Set-TextValue 'D42' '11.51'
$ws.Range('E42').Value = '  +3.79%  '
Set-TextValue 'D43' '0.6294'
$ws.Range('E43').Value = '  +4.17%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D44' '0.9998'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D45' '14.05'
$ws.Range('E45').Value = '  +3.83%  '
Set-TextValue 'D46' '0.6124'
$ws.Range('E46').Value = '  +7.19%  '
Set-TextValue 'D47' '3.780'
$ws.Range('E47').Value = '  +0.66%  '
Set-TextValue 'D48' '2.056'
$ws.Range('E48').Value = '  +3.12%  '
Set-TextValue 'D49' '125.42'
$ws.Range('E49').Value = '  +0.98%  '
Set-TextValue 'D50' '1.210'
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('E51').Value = '  +0.10%  '
